$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Olaf the Viking already occupies row 36 (Name/Gender/Role/Weapon filled in);
# fill in his stat block (Health, MP, Attack, Defense, Resistance, Skill, Speed).
$ws.Range("G36:M36").Value = 55

# Reflect where the cursor ended up after the edit.
$ws.Activate()
$ws.Range("H34").Select() | Out-Null
